# Edit script implementing the diff:
#  - update "21 years" -> "15+ years" in professional summary
#  - expand FLEEM bullet description
#  - rewrite/expand Praxis Project bullet list (4 -> 8 bullets)
#  - add a new bullet after Lake Research Partners' last bullet
#  - rewrite/expand Salsa Labs bullet list (4 -> 6 bullets)
#  - add a new bullet after Feldman Group's last bullet

$d = $word.ActiveDocument

function Replace-Text($oldText, $newText) {
    $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false,
                             $true, 1, $false, $newText, 2) | Out-Null
}

function Insert-BulletAfter($searchText, $newText) {
    $r = $d.Content
    $r.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $idx = $d.Range(0, $r.End).Paragraphs.Count
    $p = $d.Paragraphs.Item($idx)
    $p.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Item($idx + 1)
    $newPara.Range.Text = $newText
}

# ---------------------------------------------------------------------------
# 1) Professional summary: "21 years" -> "15+ years"
# ---------------------------------------------------------------------------
$old = "Distinguished Polling, Research & Redistricting Professional with 21 years of expertise"
$new = "Distinguished Polling, Research & Redistricting Professional with 15+ years of expertise"
Replace-Text $old $new

# ---------------------------------------------------------------------------
# 2) FLEEM bullet (Progressive Change Campaign Committee) expanded wording
# ---------------------------------------------------------------------------
$old = "• Conceived, architected, and engineered FLEEM web application using Twilio API for thousands of simultaneous phone calls"
$new = "• Conceived, architected, and engineered FLEEM web application using Twilio API handling tens of thousands of calls using emulated predictive dialer for regulated political surveys"
Replace-Text $old $new

# ---------------------------------------------------------------------------
# 3) Feldman Group: add a bullet right after the existing final bullet
# ---------------------------------------------------------------------------
$old = "• Enhanced value of research deliverables through advanced analytical techniques using SPSS, OSCAR, PHP, and MySQL"
$new = "• Trained staff on PHP/MySQL for data analysis and reporting systems"
Insert-BulletAfter $old $new

# ---------------------------------------------------------------------------
# 4) Praxis Project bullets: rewrite first three, expand the fourth into four
# ---------------------------------------------------------------------------
$old = "• Integrated technology solutions within organizational frameworks for social justice organizations"
$new = "• Led technology operations for multi-million dollar organization while assisting in search for full-time CTO"
Replace-Text $old $new

$old = "• Developed data management systems for community organizing efforts"
$new = "• Directed all technology decisions and practices for massive multinational non-governmental organization"
Replace-Text $old $new

$old = "• Provided technical training and support to nonprofit staff"
$new = "• Developed comprehensive frameworks for internal and external technology audits"
Replace-Text $old $new

$old = "• Built custom applications for community engagement and advocacy"
$new = "• Led training initiatives for beneficiaries on spatial and Census data analysis for public health research"
Replace-Text $old $new

$old = "• Led training initiatives for beneficiaries on spatial and Census data analysis for public health research"
$new = "• Conducted training programs for NGO staff in web development using Drupal, PHP, and MySQL"
Insert-BulletAfter $old $new

$old = "• Conducted training programs for NGO staff in web development using Drupal, PHP, and MySQL"
$new = "• Managed technology infrastructure supporting community health initiatives across multiple countries"
Insert-BulletAfter $old $new

$old = "• Managed technology infrastructure supporting community health initiatives across multiple countries"
$new = "• Architected and developed 25 Drupal sites to integrate with membership databases, activism CRMs and government agencies, under guidelines from Kellogg Foundation and Robert Wood Johnson Foundation"
Insert-BulletAfter $old $new

# ---------------------------------------------------------------------------
# 5) Lake Research Partners: add a bullet right after the existing final bullet
# ---------------------------------------------------------------------------
$old = "• Developed innovative approaches to visualizing demographic and market data for enhanced client understanding"
$new = "• Trained staff on building Python tooling for report generation and analysis"
Insert-BulletAfter $old $new

# ---------------------------------------------------------------------------
# 6) Salsa Labs bullets: rewrite first three, expand the fourth into three
# ---------------------------------------------------------------------------
$old = "• Developed software solutions for political campaigns and advocacy groups"
$new = "• Maintained and extended comprehensive geospatial analysis and reporting tools for Java-based CRM system used by tens of thousands of users simultaneously"
Replace-Text $old $new

$old = "• Built web applications for voter engagement and campaign management"
$new = "• Developed custom tile server for Web Map Service (WMS) integration using GeoTools and OpenLayers"
Replace-Text $old $new

$old = "• Integrated third-party APIs and data sources for campaign tools"
$new = "• Built advanced geospatial analysis capabilities using Java, JavaScript, MySQL, and TileMill"
Replace-Text $old $new

$old = "• Collaborated with political strategists to translate requirements into technical solutions"
$new = "• Integrated mapping and visualization tools for political campaign data analysis interfacing with Government and Activism APIs"
Replace-Text $old $new

$old = "• Integrated mapping and visualization tools for political campaign data analysis interfacing with Government and Activism APIs"
$new = "• Collaborated with political strategists to translate geospatial requirements into technical solutions"
Insert-BulletAfter $old $new

$old = "• Collaborated with political strategists to translate geospatial requirements into technical solutions"
$new = "• Handled billions of records with millions of columns in high-performance CRM system"
Insert-BulletAfter $old $new

Write-Host "Done."
